$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 4 (Argus Group Holdings Limited) - shifts subsequent rows up
$ws.Rows(4).Delete()

# Step 2: clear stale E5 (no longer present after refresh for Randall & Quilter)
$ws.Range("E5").ClearContents()

# Step 3: refresh data values for rows 2-5

# Row 2
$ws.Range("A2").Value = "Bermuda"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"
$ws.Range("C2").Value = "Insurance (General)"
$ws.Range("D2").Value = 0.2092
$ws.Range("E2").Value = -0.0707
$ws.Range("G2").Value = 0.07770246164823404
$ws.Range("H2").Value = 0.07770246164823404
$ws.Range("I2").Value = 0.02304013603203665
$ws.Range("J2").Value = 0.02225079120004012
$ws.Range("K2").Value = 6.98
$ws.Range("L2").Value = 0.004980378166250446
$ws.Range("M2").Value = 118.4167
$ws.Range("N2").Value = 0.08434237891737892
$ws.Range("O2").Value = 16.96514326647564
$ws.Range("P2").Value = 31.1167
$ws.Range("Q2").Value = 0.02216289173789174
$ws.Range("R2").Value = 4.457979942693409
$ws.Range("S2").Value = 87.30000000000001
$ws.Range("T2").Value = 0.7372270971915279
$ws.Range("U2").Value = 644.3000000000001
$ws.Range("V2").Value = 0.458903133903134
$ws.Range("W2").Value = 0.02417926003126628
$ws.Range("X2").Value = 0.05583348873153725
$ws.Range("Y2").Value = -0.03165422870027097
$ws.Range("Z2").Value = 0.7224743534333014
$ws.Range("AA2").Value = 0.01771562711477256
$ws.Range("AB2").Value = 0.04962856410047922
$ws.Range("AC2").Value = -0.03901452134495093
$ws.Range("AD2").Value = 708.5999999999999
$ws.Range("AE2").Value = 0.9612467555031695
$ws.Range("AF2").Value = 709.5612467555031
$ws.Range("AG2").Value = 65.26124675550307
$ws.Range("AH2").Value = 0.3357183274649553
$ws.Range("AI2").Value = 0.2904512904974942
$ws.Range("AJ2").Value = 0.04441772822880631
$ws.Range("AK2").Value = 0.03628323391812879
$ws.Range("AL2").Value = 51.871
$ws.Range("AM2").Value = 51.871
$ws.Range("AN2").Value = 17.64749831893009
$ws.Range("AO2").Value = 0.6207707582271405
$ws.Range("AP2").Value = 1.62531434153122
$ws.Range("AQ2").Value = 0.6207707582271405

# Row 3
$ws.Range("A3").Value = "Bermuda"
$ws.Range("B3").Value = "BF&M Limited (BER:BFM.BH)"
$ws.Range("C3").Value = "Insurance (General)"
$ws.Range("D3").Value = -0.0256
$ws.Range("E3").Value = -0.0707
$ws.Range("G3").Value = 0.1050803300043422
$ws.Range("H3").Value = 0.1050803300043422
$ws.Range("I3").Value = 0.07685627442466347
$ws.Range("J3").Value = 0.07555059112619902
$ws.Range("K3").Value = 10.9
$ws.Range("L3").Value = 0.04732957012592271
$ws.Range("M3").Value = 4.425
$ws.Range("N3").Value = 0.025
$ws.Range("O3").Value = 0.4059633027522935
$ws.Range("P3").Value = 4.425
$ws.Range("Q3").Value = 0.025
$ws.Range("R3").Value = 0.4059633027522935
$ws.Range("S3").Value = 0.0
$ws.Range("T3").Value = 0.0
$ws.Range("U3").Value = 119.2
$ws.Range("V3").Value = 0.6734463276836158
$ws.Range("W3").Value = 0.03828591499824377
$ws.Range("X3").Value = 0.04592464886435726
$ws.Range("Y3").Value = -0.007638733866113495
$ws.Range("Z3").Value = 1.629865534324133
$ws.Range("AA3").Value = 0.1231373045744065
$ws.Range("AB3").Value = 0.04592464886435726
$ws.Range("AC3").Value = 0.07721265571004923
$ws.Range("AD3").Value = 0.0
$ws.Range("AE3").Value = 0.0
$ws.Range("AF3").Value = 0.0
$ws.Range("AG3").Value = -119.2
$ws.Range("AH3").Value = 0.0
$ws.Range("AI3").Value = 0.0
$ws.Range("AJ3").Value = -2.062283737024222
$ws.Range("AK3").Value = -0.5766811804547654
$ws.Range("AL3").Value = 0.071
$ws.Range("AM3").Value = 0.071
$ws.Range("AN3").Value = 0.0
$ws.Range("AO3").Value = 249.2957746478873
$ws.Range("AP3").Value = -6.273684210526316
$ws.Range("AQ3").Value = 249.2957746478873

# Row 4
$ws.Range("A4").Value = "Bermuda"
$ws.Range("B4").Value = "Watford Holdings Ltd. (NasdaqGS:WTRE)"
$ws.Range("C4").Value = "Insurance (General)"
$ws.Range("G4").Value = 0.1624902723735409
$ws.Range("H4").Value = 0.1624902723735409
$ws.Range("I4").Value = 0.04481050684653598
$ws.Range("J4").Value = 0.04481050684653598
$ws.Range("K4").Value = -13.2
$ws.Range("L4").Value = -0.02054474708171206
$ws.Range("M4").Value = 77.9
$ws.Range("N4").Value = 0.1132102892021509
$ws.Range("O4").Value = -5.901515151515152
$ws.Range("P4").Value = -0.0
$ws.Range("Q4").Value = -0.0
$ws.Range("R4").Value = 0.0
$ws.Range("S4").Value = 77.9
$ws.Range("T4").Value = 1.0
$ws.Range("U4").Value = 195.3
$ws.Range("V4").Value = 0.2838250254323499
$ws.Range("W4").Value = -0.01373855120732723
$ws.Range("X4").Value = 0.0759312724739381
$ws.Range("Y4").Value = -0.08966982368126533
$ws.Range("Z4").Value = 0.3953453857472278
$ws.Range("AA4").Value = 0.01771562711477256
$ws.Range("AB4").Value = 0.05673014845972348
$ws.Range("AC4").Value = -0.03901452134495093
$ws.Range("AD4").Value = 562.8
$ws.Range("AE4").Value = 0.9612467555031695
$ws.Range("AF4").Value = 563.7612467555032
$ws.Range("AG4").Value = 368.4612467555032
$ws.Range("AH4").Value = 0.4503384446292469
$ws.Range("AI4").Value = 0.3801334894218598
$ws.Range("AJ4").Value = 0.3487362875431755
$ws.Range("AK4").Value = 0.2861254348846389
$ws.Range("AL4").Value = 39.4
$ws.Range("AM4").Value = 39.4
$ws.Range("AN4").Value = 19.41827968119242
$ws.Range("AO4").Value = 0.7284263959390863
$ws.Range("AP4").Value = 12.71301268866243
$ws.Range("AQ4").Value = 0.7284263959390863

# Row 5
$ws.Range("A5").Value = "Bermuda"
$ws.Range("B5").Value = "Randall & Quilter Investment Holdings Ltd (AIM:RQIH)"
$ws.Range("C5").Value = "Insurance (General)"
$ws.Range("D5").Value = 0.444
$ws.Range("G5").Value = -0.03726120673349725
$ws.Range("H5").Value = -0.03726120673349725
$ws.Range("I5").Value = -0.02685833175714015
$ws.Range("J5").Value = -0.02455415462124128
$ws.Range("K5").Value = 9.28
$ws.Range("L5").Value = 0.01755248723283525
$ws.Range("M5").Value = 36.0917
$ws.Range("N5").Value = 0.06697290777509743
$ws.Range("O5").Value = 3.889191810344828
$ws.Range("P5").Value = 26.6917
$ws.Range("Q5").Value = 0.04952996845425868
$ws.Range("R5").Value = 2.876260775862069
$ws.Range("S5").Value = 9.400000000000002
$ws.Range("T5").Value = 0.2604476929598772
$ws.Range("U5").Value = 329.8
$ws.Range("V5").Value = 0.6119873817034701
$ws.Range("W5").Value = 0.02417926003126628
$ws.Range("X5").Value = 0.05583348873153725
$ws.Range("Y5").Value = -0.03165422870027097
$ws.Range("Z5").Value = 3.049019607843138
$ws.Range("AA5").Value = -0.07486609889417685
$ws.Range("AB5").Value = 0.04962856410047922
$ws.Range("AC5").Value = -0.1244946629946561
$ws.Range("AD5").Value = 145.8
$ws.Range("AE5").Value = 0.0
$ws.Range("AF5").Value = 145.8
$ws.Range("AG5").Value = -184.0
$ws.Range("AH5").Value = 0.2129399737111144
$ws.Range("AI5").Value = 0.2299684542586751
$ws.Range("AJ5").Value = -0.5184559030712877
$ws.Range("AK5").Value = -0.6048652202498357
$ws.Range("AL5").Value = 12.4
$ws.Range("AM5").Value = 12.4
$ws.Range("AN5").Value = -18.62068965517242
$ws.Range("AO5").Value = -1.145161290322581
$ws.Range("AP5").Value = 23.49936143039591
$ws.Range("AQ5").Value = -1.145161290322581
